# Investment ticker symbols workbook: add a "Group" column to Stock_list,
# add header rows to Index_list / Commodity_list, and wire up the
# dynamic dropdown UI (headers used for data-validation lists).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Stock_list: insert a header row, add a Group column (C)
# ---------------------------------------------------------------------
$stock = $wb.Worksheets.Item("Stock_list")

# Shift existing data down by one row and add the new header row.
$stock.Rows.Item(1).Insert() | Out-Null
$stock.Range("A1").Value = "Ticker"
$stock.Range("B1").Value = "Stock"
$stock.Range("C1").Value = "Group"

$groups = @(
    'Chemical & Functional Materials',
    'Chemical & Functional Materials',
    'Retailers',
    'Pharmaceuticals',
    'Consumer & Medical Goods',
    'Chemical & Functional Materials',
    'Pharmaceuticals',
    'Consumer & Medical Goods',
    'Oil, Gas, & Energy',
    'Oil, Gas, & Energy',
    'Pharmaceuticals',
    'Oil, Gas, & Energy',
    'Oil, Gas, & Energy',
    'Oil, Gas, & Energy',
    'Consumer & Medical Goods',
    'Consumer & Medical Goods',
    'Oil, Gas, & Energy',
    'Chemical & Functional Materials',
    'Chemical & Functional Materials',
    'Chemical & Functional Materials',
    'Chemical & Functional Materials',
    'Chemical & Functional Materials',
    'Oil, Gas, & Energy',
    'Chemical & Functional Materials',
    'Food & Fragrances',
    'Pharmaceuticals',
    'Chemical & Functional Materials',
    'Food & Fragrances',
    'Food & Fragrances',
    'Chemical & Functional Materials',
    'Chemical & Functional Materials',
    'Cosmetics',
    'Chemical & Functional Materials',
    'Oil, Gas, & Energy',
    'Pharmaceuticals',
    'Chemical & Functional Materials',
    'Food & Fragrances',
    'Food & Fragrances',
    'Oil, Gas, & Energy',
    'Oil, Gas, & Energy',
    'Pharmaceuticals',
    'Consumer & Medical Goods',
    'Pharmaceuticals',
    'Oil, Gas, & Energy',
    'Oil, Gas, & Energy',
    'Cosmetics',
    'Chemical & Functional Materials',
    'Retailers',
    'Consumer & Medical Goods',
    'Chemical & Functional Materials',
    'Consumer & Medical Goods',
    'Retailers'
)

for ($i = 0; $i -lt $groups.Length; $i++) {
    $stock.Cells.Item(2 + $i, 3).Value = $groups[$i]
}

# Column widths for the new Stock / Group columns (best-fit widths, as
# computed by Excel's own AutoFit for this content/font).
$stock.Columns.Item(2).ColumnWidth = 61.17
$stock.Columns.Item(3).ColumnWidth = 25

# ---------------------------------------------------------------------
# Index_list: insert a header row
# ---------------------------------------------------------------------
$index = $wb.Worksheets.Item("Index_list")
$index.Rows.Item(1).Insert() | Out-Null
$index.Range("A1").Value = "Ticker"
$index.Range("B1").Value = "Index"
$index.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------
# Commodity_list: add a header row
# ---------------------------------------------------------------------
$commodity = $wb.Worksheets.Item("Commodity_list")
$commodity.Range("A1").Value = "Ticker"
$commodity.Range("B1").Value = "Commodity"
$commodity.Range("D7").Select() | Out-Null

# ---------------------------------------------------------------------
# Selection / active sheet state
# ---------------------------------------------------------------------
$stock.Activate() | Out-Null
$stock.Range("C53").Select() | Out-Null
